# The underlying edit swaps the data held in row 2 and row 3 of the sheet
# (title, timestamp, historical distance) while the time bucket ("day_31_beyond")
# is identical for both rows already. The uri column's displayed text is swapped
# too, though the underlying hyperlink targets are left exactly as they were.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current row 2 and row 3 values (use Value2 - Value is unreliable here)
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$e2 = $ws.Range("E2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$e3 = $ws.Range("E3").Value2

# Row 2 gets what used to be row 3's content
$ws.Range("A2").Value2 = $a3
$ws.Range("B2").Value2 = $b3
$ws.Range("C2").Value2 = $c3
$ws.Range("E2").Value2 = $e3

# Row 3 gets what used to be row 2's content
$ws.Range("A3").Value2 = $a2
$ws.Range("B3").Value2 = $b2
$ws.Range("C3").Value2 = $c2
$ws.Range("E3").Value2 = $e2
